$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 3's Mode value first (C2/"serial" stays put; this just appends
# "parallel" as a brand new shared string).
$ws.Range("C3").Value = "parallel"

# Row 2: rename the existing REST test to the new Weather REST API test name
$ws.Range("A2").Value = "WeatherRESTAPITest"

# Apply the new "small Segoe UI, vertically centered" look to A2
$ws.Range("A2").Font.Size = 9
$ws.Range("A2").Font.Name = "Segoe UI"
$ws.Range("A2").VerticalAlignment = -4108

# Match A2's formatting on the new row's test-name cell (copy/paste the
# formatting instead of re-applying each font property so the two cells
# end up sharing the very same cell style)
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New row 3: EnergyUITest / YES / (parallel already set above)
$ws.Range("A3").Value = "EnergyUITest"
$ws.Range("B3").Value = "YES"

# Move the active selection down to the next empty row, like a user who just
# finished typing row 3 and tabbed/entered down to A4
$ws.Range("A4").Select()

Write-Output "done"
